$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF2").Value = "Digit_before_and_AFACT"
$ws.Range("AF3").Value = "MAB_phase"
$ws.Range("AF4").Value = "dichotic_phase"
$ws.Range("AF5").Value = "dichotic_phase"
$ws.Range("AF6").Value = "Dichotic_and_AFACT"
$ws.Range("AF7").Value = "Dichotic_and_AFACT"
$ws.Range("AF8").Value = "MAB_and_Digit_after"
$ws.Range("AF9").Value = "Digit_before_and_AFACT"
$ws.Range("AF10").Value = "Digit_before_and_AFACT"
$ws.Range("AF11").Value = "Dichotic_and_AFACT"
$ws.Range("AF13").Value = "Dichotic_and_AFACT"
$ws.Range("AF14").Value = "MAB_and_AFACT"
$ws.Range("AF15").Value = "Dichotic_and_AFACT"
$ws.Range("AF16").Value = "dichotic_phase"
$ws.Range("AF17").Value = "MAB_and_Digit_after"
$ws.Range("AF19").Value = "Digit_before_and_AFACT"
$ws.Range("AF20").Value = "Dichotic_and_AFACT"
$ws.Range("AF22").Value = "dichotic_phase"
$ws.Range("AF23").Value = "Dichotic_and_AFACT"
$ws.Range("AF24").Value = "MAB_and_AFACT"
$ws.Range("AF25").Value = "MAB_phase"
$ws.Range("AF26").Value = "dichotic_phase"
$ws.Range("AF27").Value = "Digit_before_and_AFACT"
$ws.Range("AF28").Value = "dichotic_phase"
$ws.Range("AF32").Value = "dichotic_phase"
$ws.Range("AF33").Value = "Dichotic_and_AFACT"
$ws.Range("AF34").Value = "MAB_phase"
$ws.Range("AF35").Value = "MAB_and_AFACT"
$ws.Range("AF36").Value = "Dichotic_and_AFACT"
$ws.Range("AF37").Value = "MAB_and_AFACT"
$ws.Range("AF38").Value = "dichotic_phase"
$ws.Range("AF40").Value = "MAB_and_Digit_after"
$ws.Range("AF41").Value = "MAB_and_Digit_after"
$ws.Range("AF42").Value = "Dichotic_and_AFACT"
$ws.Range("AF43").Value = "MAB_and_Digit_after"
$ws.Range("AF45").Value = "dichotic_phase"
$ws.Range("AF46").Value = "MAB_and_Digit_after"
$ws.Range("AF47").Value = "MAB_phase"
$ws.Range("AF48").Value = "MAB_phase"
$ws.Range("AF49").Value = "Digit_before_and_AFACT"
$ws.Range("AF51").Value = "Dichotic_and_AFACT"
$ws.Range("AF52").Value = "Dichotic_and_AFACT"
$ws.Range("AF53").Value = "Dichotic_and_AFACT"
$ws.Range("AF54").Value = "Digit_before_and_AFACT"
$ws.Range("AF55").Value = "MAB_and_Digit_after"
$ws.Range("AF56").Value = "dichotic_phase"
$ws.Range("AF57").Value = "MAB_phase"
$ws.Range("AF59").Value = "dichotic_phase"
$ws.Range("AF61").Value = "MAB_and_AFACT"
$ws.Range("AF62").Value = "dichotic_phase"
$ws.Range("AF63").Value = "dichotic_phase"
$ws.Range("AF64").Value = "Dichotic_and_AFACT"
$ws.Range("AF66").Value = "MAB_and_AFACT"
$ws.Range("AF67").Value = "Dichotic_and_AFACT"
$ws.Range("AF68").Value = "Digit_before_and_AFACT"
$ws.Range("AF70").Value = "dichotic_phase"
$ws.Range("AF71").Value = "Dichotic_and_AFACT"
$ws.Range("AF72").Value = "MAB_phase"
$ws.Range("AF73").Value = "dichotic_phase"
$ws.Range("AF74").Value = "MAB_and_Digit_after"
$ws.Range("AF75").Value = "MAB_and_AFACT"
$ws.Range("AF76").Value = "Digit_before_and_AFACT"
$ws.Range("AF78").Value = "dichotic_phase"
$ws.Range("AF79").Value = "Dichotic_and_AFACT"
$ws.Range("AF80").Value = "MAB_and_Digit_after"
$ws.Range("AF81").Value = "MAB_phase"
